$wb = $excel.ActiveWorkbook

# Update Status text ("Ready for handoff" -> "Handback transform failed") for the
# 68d5ce78-... row (row 3) on every sheet that references it (Overview, zh-cn, de-de).
foreach ($sheetName in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C3").Value = "Handback transform failed"
}

# Add the new "Error Detail" (column L) messages produced by the handback transform
# for row 3 (the 68d5ce78-... file) on the zh-cn and de-de sheets.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("L3").Value = "Handback file name: i2zwil2w.vxc is different with handoff file name: 68d5ce78-33d4-4845-a5ff-6a3d8e201d14.16f5bb49eb9311ac9dc863c3b1553b871ff1c912.zh-cn."

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("L3").Value = "Handback file name: i2zwil2w.vxc is different with handoff file name: 68d5ce78-33d4-4845-a5ff-6a3d8e201d14.16f5bb49eb9311ac9dc863c3b1553b871ff1c912.de-de."
